# "Generate Report for Handoff"
#
# The localization CI job re-ran and produced a fresh report:
#   - status flips from "In Translation" to "Ready for handoff"
#     (Overview!E2/F2, zh-cn!C2, de-de!C2 all share this string)
#   - the two "Latest ... Datetime" timestamps tick forward a minute
#     (Overview!G2, zh-cn!H2, de-de!H2)
#   - the two "...Datetime" columns were widened slightly in all three
#     sheets to fit the new label

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 15:12:45"

# zh-cn / de-de status + handoff-datetime columns widened (col E & F)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 15:12:37"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 15:12:45"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
